$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that are being updated, so purely
# numeric-looking values (e.g. "10.10", "0.999") keep their exact text
# representation rather than being auto-converted to numbers by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D42", "D43", "D45", "D48", "D49", "D50")
foreach ($pc in $priceCells) {
    $ws.Range($pc).NumberFormat = "@"
}

# Apply the updated coin data.
$ws.Range("D2").Value = "58.477.11"
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("D3").Value = "2.537.26"
$ws.Range("E3").Value = "  -3.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "507.11"
$ws.Range("E5").Value = "  -4.26%  "
$ws.Range("D6").Value = "143.99"
$ws.Range("E6").Value = "  -7.38%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -4.48%  "
$ws.Range("D9").Value = "2.541.83"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("D10").Value = "6.08"
$ws.Range("E10").Value = "  -8.95%  "
$ws.Range("E11").Value = "  -6.86%  "
$ws.Range("E12").Value = "  -5.62%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "2.982.31"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "58.459.91"
$ws.Range("E15").Value = "  -4.16%  "
$ws.Range("D16").Value = "20.67"
$ws.Range("E16").Value = "  -6.05%  "
$ws.Range("E17").Value = "  -6.49%  "
$ws.Range("D18").Value = "2.536.47"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").Value = "334.58"
$ws.Range("E20").Value = "  -5.51%  "
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  -5.07%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("D24").Value = "60.58"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "0.410"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -5.41%  "
$ws.Range("D28").Value = "2.649.41"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").Value = "0.0₃0786"
$ws.Range("E29").Value = "  -9.49%  "
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "149.72"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "5.84"
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.53"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D36").Value = "0.915"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("E37").Value = "  -6.33%  "
$ws.Range("E38").Value = "  -7.71%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -12.08%  "
$ws.Range("E41").Value = "  -7.19%  "
$ws.Range("D42").Value = "284.08"
$ws.Range("E42").Value = "  -7.36%  "
$ws.Range("D43").Value = "3.53"
$ws.Range("E43").Value = "  -8.03%  "
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("E47").Value = "  -5.41%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "18.64"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.29"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "0.0226"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("E51").Value = "  -10.01%  "
